# F224A -> F224I revision:
#   - cover title and the board-id table cell get their product code
#     updated from F224A to F224I
#   - the language cell changes from German to Italian
#   - the previously empty "Notes" line gets a note added

$d = $word.ActiveDocument

# 1) Big cover-page title. "F224A" there is split across three separate
#    runs ("F" / "224" / "A"); only the trailing "A" run actually changes
#    (to "I"). Scope the Find to the title paragraph so only that run's
#    text is touched and the existing "F"/"224" run boundaries are left
#    completely alone.
$titleRange = $d.Paragraphs.Item(1).Range
$titleRange.Find.Execute("A", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "I", 2) | Out-Null

# 2) Table cell with the board/product code: "F224A" -> "F224I".
$d.Content.Find.Execute("F224A", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "F224I", 2) | Out-Null

# 3) Table cell with the language name: "German" -> "Italian".
$d.Content.Find.Execute("German", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Italian", 2) | Out-Null

# 4) The blank paragraph right after "Notes: " gets filled in. Locate the
#    "Notes: " paragraph, then its following (empty) paragraph, and set
#    that paragraph's text directly -- trimming the trailing paragraph
#    mark off the range first so no new paragraph gets inserted and the
#    existing (empty) run's formatting is kept.
$notesLabel = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Notes: `r") {
        $notesLabel = $p
        break
    }
}
if ($notesLabel -ne $null) {
    $notesBody = $notesLabel.Next()
    $bodyRange = $notesBody.Range.Duplicate
    $bodyRange.End = $bodyRange.End - 1
    $bodyRange.Text = "scratches on top"
}
